# ajustes sanity semilla 6 en clases de portabilidad prepago y postpago
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepago / Postpago portability class values were rotated to new test data.
# Write order matches the order these new shared-string values were first
# introduced in the saved workbook.
$ws.Range("H9").Value  = "662496115"
$ws.Range("E10").Value = "3046008600"
$ws.Range("I2").Value  = "http://10.69.60.74:8080/PortabilidadServiceEAR-HPNPCommunicationsDelegateEJB/NPCRMWSImpl"
$ws.Range("B14").Value = "662496115"

# Reflect the cell the user was last working in (I2) as the active selection.
$ws.Range("I2").Select()

# Best-effort: scroll the view so column C is left-most (matches author's
# saved view state). Not all hosts persist window/scroll geometry to xlsx.
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollColumn = 3
    $win.ScrollRow = 1
}
